$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster-pair rows (rows 2-17) reflecting the updated TPM data and the
# newly introduced "Resolving-Mac" cluster. Columns:
# A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T numeric metrics.
$data = @(
    @("ECs","Sema3a","Plxna1","ECs",3,1,0.8200746666666667,2.460224,0.04311293902675128,0.04311293902675128,3,1,7.368211,22.104633,0.1280150446959744,0.1280150446959744,6.042483179754667,54.382348617792,0.005519104816484387,0.005519104816484387),
    @("ECs","Sema3a","Plxna1","FAPs",3,1,0.8200746666666667,2.460224,0.04311293902675128,0.04311293902675128,3,1,18.05628333333334,54.16885000000001,0.3137092461059877,0.3137092461059876,14.80750053582223,133.2675048224,0.01352492759949556,0.01352492759949556),
    @("ECs","Sema3a","Plxna1","MuSCs",3,1,0.8200746666666667,2.460224,0.04311293902675128,0.04311293902675128,3,1,18.794071,56.38221299999999,0.3265275436716344,0.3265275436716343,15.41254151063467,138.712873595712,0.01407756208087004,0.01407756208087004),
    @("ECs","Sema3a","Plxna1","Resolving-Mac",3,1,0.8200746666666667,2.460224,0.04311293902675128,0.04311293902675128,3,1,13.338818,40.016454,0.2317481655264036,0.2317481655264036,10.93882672507733,98.44944052569599,0.009991344529901304,0.009991344529901304),
    @("FAPs","Sema3a","Plxna1","ECs",3,1,1.358031333333334,4.074094000000001,0.07139437962203982,0.07139437962203982,3,1,7.368211,22.104633,0.1280150446959744,0.1280150446959744,10.00626140861134,90.056352677502,0.009139554698356794,0.009139554698356794),
    @("FAPs","Sema3a","Plxna1","FAPs",3,1,1.358031333333334,4.074094000000001,0.07139437962203982,0.07139437962203982,3,1,18.05628333333334,54.16885000000001,0.3137092461059877,0.3137092461059876,24.52099853021112,220.6889867719001,0.0223970770074348,0.02239707700743479),
    @("FAPs","Sema3a","Plxna1","MuSCs",3,1,1.358031333333334,4.074094000000001,0.07139437962203982,0.07139437962203982,3,1,18.794071,56.38221299999999,0.3265275436716344,0.3265275436716343,25.52293729889134,229.706435690022,0.02331223140994485,0.02331223140994485),
    @("FAPs","Sema3a","Plxna1","Resolving-Mac",3,1,1.358031333333334,4.074094000000001,0.07139437962203982,0.07139437962203982,3,1,13.338818,40.016454,0.2317481655264036,0.2317481655264036,18.11453279363067,163.030795142676,0.01654551650630338,0.01654551650630338),
    @("MuSCs","Sema3a","Plxna1","ECs",3,1,16.8273,50.4819,0.8846442749337277,0.8846442749337278,3,1,7.368211,22.104633,0.1280150446959744,0.1280150446959744,123.9870969603,1115.8838726427,0.1132477763956791,0.1132477763956791),
    @("MuSCs","Sema3a","Plxna1","FAPs",3,1,16.8273,50.4819,0.8846442749337277,0.8846442749337278,3,1,18.05628333333334,54.16885000000001,0.3137092461059877,0.3137092461059876,303.838496535,2734.546468815,0.2775210885614378,0.2775210885614378),
    @("MuSCs","Sema3a","Plxna1","MuSCs",3,1,16.8273,50.4819,0.8846442749337277,0.8846442749337278,3,1,18.794071,56.38221299999999,0.3265275436716344,0.3265275436716343,316.2534709382999,2846.281238444699,0.2888607221172841,0.2888607221172841),
    @("MuSCs","Sema3a","Plxna1","Resolving-Mac",3,1,16.8273,50.4819,0.8846442749337277,0.8846442749337278,3,1,13.338818,40.016454,0.2317481655264036,0.2317481655264036,224.4562921313999,2020.1066291826,0.2050146878593269,0.2050146878593269),
    @("Resolving-Mac","Sema3a","Plxna1","ECs",1,0.3333333333333333,0.016138,0.048414,0.0008484064174811467,0.0008484064174811467,3,1,7.368211,22.104633,0.1280150446959744,0.1280150446959744,0.118908189118,1.070173702062,0.0001086087854542005,0.0001086087854542005),
    @("Resolving-Mac","Sema3a","Plxna1","FAPs",1,0.3333333333333333,0.016138,0.048414,0.0008484064174811467,0.0008484064174811467,3,1,18.05628333333334,54.16885000000001,0.3137092461059877,0.3137092461059876,0.2913923004333334,2.6225307039,0.0002661529376194924,0.0002661529376194923),
    @("Resolving-Mac","Sema3a","Plxna1","MuSCs",1,0.3333333333333333,0.016138,0.048414,0.0008484064174811467,0.0008484064174811467,3,1,18.794071,56.38221299999999,0.3265275436716344,0.3265275436716343,0.303298717798,2.729688460182,0.00027702806353537,0.00027702806353537),
    @("Resolving-Mac","Sema3a","Plxna1","Resolving-Mac",1,0.3333333333333333,0.016138,0.048414,0.0008484064174811467,0.0008484064174811467,3,1,13.338818,40.016454,0.2317481655264036,0.2317481655264036,0.215261844884,1.937356603956,0.0001966166308720839,0.0001966166308720839)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
